# Applies the "Correction type pour génération à partir fsh" edit:
#   - Metadata sheet, "Name" row (row 4): set the value cell (B4) to
#     "EnsemblesavoirfaireCisisVs" (was empty).
#   - Metadata sheet, "Date" row (row 8): update the value cell (B8) from
#     "2025-07-17T14:35:50+00:00" to "2025-07-18T06:40:38+00:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "EnsemblesavoirfaireCisisVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
